$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Austin")

$ws.Range("A23").Value = "4/4/20"
$ws.Range("B23").Value = 460
$ws.Range("C23").Value = "Travis County"

$ws.Range("A24").Select()
